# Apply: shift annual cashflow dates down a year on "PBO" and "Service Cost"
# sheets (matching the already-correct "PVFB" sheet), dropping the
# EOMONTH() formula in the final row in favor of a plain static value.

$wb = $excel.ActiveWorkbook

# Build the target date sequence: 2021-12-31 ... 2100-12-31 (80 year-end
# dates), taken straight from the reference "PVFB" sheet which already has
# the desired values.
$refWs = $wb.Worksheets.Item("PVFB")

$sheetNames = @("PBO", "Service Cost")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    for ($r = 2; $r -le 81; $r++) {
        $ws.Cells.Item($r, 1).Value = $refWs.Cells.Item($r, 1).Value()
    }

    $ws.Range("A2:A81").Select()
}

# Restore the view state / selection to match the target workbook.
$pbo = $wb.Worksheets.Item("PBO")
$pbo.Activate()
$pbo.Range("A2").Select()

$sc = $wb.Worksheets.Item("Service Cost")
$sc.Activate()
$sc.Application.ActiveWindow.ScrollRow = 15
$sc.Range("A2").Select()

$pbo.Activate()
